$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'38.93"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.04%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.098"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.78%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08081"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.931"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.58%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.175"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'7.968"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.28%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9307"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'0.31%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1479"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.36%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1929"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.18%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09124"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.33%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03500"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.33%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09783"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.31%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001391"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.02%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005857"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.26%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.789"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.40%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D19").Value = "'0.3425"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.42%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1303"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'0.76%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.755"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.21%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'3.12%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-0.41%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001236"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.11%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004285"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-12.89%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'0.01%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D39").Value = "'0.02044"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.79%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05085"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.47%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007448"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.73%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.01027"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.83%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1350"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-1.60%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002121"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-0.47%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009115"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-6.32%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006197"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.20%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003102"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Value = "'0.00002101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002001"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.04%"
$ws.Range("E51").Style = "Normal"
